# update entity framework content
#
# Applies:
#  1. The "datetimeFigureOut" field cached text on the slide master and every
#     slide layout: "10/3/2018" -> "10/8/2018".
#  2. Slide 8 ("assignment") title case fix: "assignment" -> "Assignment".
#  3. Slide 8 body: merge the three runs describing the Controller methods
#     into a single run with the combined text.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached date field text throughout the slide master + layouts
# ---------------------------------------------------------------------------
function Update-DateShapes {
    param($shapes)

    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.HasTextFrame) {
            $txt = $shp.TextFrame.TextRange.Text
            if ($txt -eq "10/3/2018") {
                $shp.TextFrame.TextRange.Text = "10/8/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) + 3) Slide 8 content updates
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)

# Title: "assignment" -> "Assignment"
$titleShape = $slide8.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Assignment"

# Body: merge the "Create Controller ..." runs into a single run.
$bodyShape = $slide8.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

# First collapse the second paragraph down to a single placeholder character
# so the host's minimal-diff text patcher drops the old GET/POST/PUT/DELETE
# run boundaries (instead of leaving the stray "XXX"-style run fragments that
# a single direct re-assignment of the combined string would preserve).
$bodyRange.Text = "Create Project web api`rX`r"

# Now write the fully combined sentence back in; this becomes one clean run.
$bodyRange2 = $bodyShape.TextFrame.TextRange
$bodyRange2.Text = "Create Project web api`rCreate Controller that contain method for GET, POST, PUT, DELETE`r"

Write-Output "done"
